$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 750
$ws.Range("I40").Value = 750
$ws.Range("K40").Value = 750
$ws.Range("M40").Value = -575

$ws.Range("H48").Value = 4875
$ws.Range("J48").Value = 5500
$ws.Range("L48").Value = 16500
$ws.Range("N48").Value = -17084

$ws.Range("H56").Value = 4875
$ws.Range("J56").Value = 5500
$ws.Range("L56").Value = 16500
$ws.Range("N56").Value = -17568

$ws.Range("H64").Value = 4284
$ws.Range("J64").Value = 3960
$ws.Range("L64").Value = 3960
$ws.Range("N64").Value = -4456

$ws.Range("H67").Value = 4284
$ws.Range("J67").Value = 3960
$ws.Range("L67").Value = 3960
$ws.Range("N67").Value = -5676

$ws.Range("H96").Value = 1939.9
$ws.Range("I96").Value = 2665
$ws.Range("K96").Value = 7995
$ws.Range("M96").Value = -6622

$ws.Range("H100").Value = 2253.3333
$ws.Range("J100").Value = 2896
$ws.Range("L100").Value = 2896
$ws.Range("N100").Value = -3978

$ws.Range("H132").Value = 8555681
$ws.Range("I132").Value = 10422527
$ws.Range("K132").Value = 31267581
$ws.Range("M132").Value = -31265051

$ws.Range("H137").Value = 2527.475
$ws.Range("I137").Value = 1870.7646
$ws.Range("J137").Value = 3012.8696
$ws.Range("K137").Value = 5612.293799999999
$ws.Range("L137").Value = 9038.6088
$ws.Range("M137").Value = -3062.293799999999
$ws.Range("N137").Value = -14138.6088

$ws.Range("H138").Value = 2498.9136
$ws.Range("I138").Value = 1775.4375
$ws.Range("J138").Value = 2677
$ws.Range("K138").Value = 5326.3125
$ws.Range("L138").Value = 8031
$ws.Range("M138").Value = -186.3125
$ws.Range("N138").Value = -18311

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1272.7273
$ws.Range("I2").Value = 1075
$ws.Range("J2").Value = 1800
$ws.Range("K2").Value = 1075
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = -962
$ws.Range("N2").Value = -2026

$ws.Range("H32").Value = 4683.05
$ws.Range("I32").Value = 4002.1099
$ws.Range("J32").Value = 11568.111
$ws.Range("K32").Value = 4002.1099
$ws.Range("L32").Value = 11568.111
$ws.Range("M32").Value = -3715.1099
$ws.Range("N32").Value = -12142.111

$ws.Range("H61").Value = 2204.75
$ws.Range("I61").Value = 1942.2
$ws.Range("J61").Value = 2642.3333
$ws.Range("K61").Value = 1942.2
$ws.Range("L61").Value = 2642.3333
$ws.Range("M61").Value = -1730.2
$ws.Range("N61").Value = -3066.3333

$ws.Range("H62").Value = 55960
$ws.Range("J62").Value = 55960
$ws.Range("L62").Value = 55960
$ws.Range("N62").Value = -57208

$ws.Range("H65").Value = 55960
$ws.Range("J65").Value = 55960
$ws.Range("L65").Value = 167880
$ws.Range("N65").Value = -174120

$ws.Range("H74").Value = 1695.3334
$ws.Range("I74").Value = 889.6667
$ws.Range("K74").Value = 889.6667
$ws.Range("M74").Value = -15.66669999999999

$ws.Range("H77").Value = 1695.3334
$ws.Range("I77").Value = 889.6667
$ws.Range("K77").Value = 4448.3335
$ws.Range("M77").Value = -80.33349999999973

$ws.Range("H88").Value = 2000
$ws.Range("J88").Value = 2000
$ws.Range("L88").Value = 2000
$ws.Range("N88").Value = -2812

$ws.Range("H91").Value = 2000
$ws.Range("J91").Value = 2000
$ws.Range("L91").Value = 2000
$ws.Range("N91").Value = -4808

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H116").Value = 1272.7273
$ws.Range("I116").Value = 1075
$ws.Range("J116").Value = 1800
$ws.Range("K116").Value = 1075
$ws.Range("L116").Value = 1800
$ws.Range("M116").Value = 1219
$ws.Range("N116").Value = -6388

$ws.Range("H132").Value = 4213.625
$ws.Range("I132").Value = 4300.067
$ws.Range("J132").Value = 4069.5557
$ws.Range("K132").Value = 12900.201
$ws.Range("L132").Value = 12208.6671
$ws.Range("M132").Value = -10370.201
$ws.Range("N132").Value = -17268.6671

$ws.Range("H136").Value = 2204.75
$ws.Range("I136").Value = 1942.2
$ws.Range("J136").Value = 2642.3333
$ws.Range("K136").Value = 5826.6
$ws.Range("L136").Value = 7926.999899999999
$ws.Range("M136").Value = -3276.6
$ws.Range("N136").Value = -13026.9999

$ws.Range("H139").Value = 36702.145
$ws.Range("J139").Value = 36702.145
$ws.Range("L139").Value = 36702.145
$ws.Range("N139").Value = -46982.145

$ws.Range("H140").Value = 36244.75
$ws.Range("J140").Value = 36244.75
$ws.Range("L140").Value = 36244.75
$ws.Range("N140").Value = -46604.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1272.7273
$ws.Range("I3").Value = 1075
$ws.Range("J3").Value = 1800
$ws.Range("K3").Value = 1075
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = -961
$ws.Range("N3").Value = -2028

$ws.Range("H86").Value = 3905.3635
$ws.Range("I86").Value = 4044.0625
$ws.Range("J86").Value = 3535.5
$ws.Range("K86").Value = 4044.0625
$ws.Range("L86").Value = 3535.5
$ws.Range("M86").Value = -2921.0625
$ws.Range("N86").Value = -5781.5

$ws.Range("H89").Value = 3905.3635
$ws.Range("I89").Value = 4044.0625
$ws.Range("J89").Value = 3535.5
$ws.Range("K89").Value = 20220.3125
$ws.Range("L89").Value = 17677.5
$ws.Range("M89").Value = -14604.3125
$ws.Range("N89").Value = -28909.5

$ws.Range("H134").Value = 5171.96
$ws.Range("I134").Value = 1130.0526
$ws.Range("J134").Value = 17971.334
$ws.Range("K134").Value = 3390.1578
$ws.Range("L134").Value = 53914.00199999999
$ws.Range("M134").Value = -855.1578
$ws.Range("N134").Value = -58984.00199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1534.3019
$ws.Range("I31").Value = 1344.3556
$ws.Range("J31").Value = 2602.75
$ws.Range("K31").Value = 1344.3556
$ws.Range("L31").Value = 2602.75
$ws.Range("M31").Value = -1049.3556
$ws.Range("N31").Value = -3192.75

$ws.Range("H34").Value = 1534.3019
$ws.Range("I34").Value = 1344.3556
$ws.Range("J34").Value = 2602.75
$ws.Range("K34").Value = 1344.3556
$ws.Range("L34").Value = 2602.75
$ws.Range("M34").Value = -1142.3556
$ws.Range("N34").Value = -3006.75

$ws.Range("H106").Value = 27062.25
$ws.Range("J106").Value = 27062.25
$ws.Range("L106").Value = 27062.25
$ws.Range("N106").Value = -29586.25

$ws.Range("H107").Value = 753.1786
$ws.Range("I107").Value = 451.81818
$ws.Range("K107").Value = 451.81818
$ws.Range("M107").Value = 1468.18182

$ws.Range("H132").Value = 1677.6552
$ws.Range("I132").Value = 1364.6471
$ws.Range("J132").Value = 2121.0833
$ws.Range("K132").Value = 4093.9413
$ws.Range("L132").Value = 6363.249899999999
$ws.Range("M132").Value = -1563.9413
$ws.Range("N132").Value = -11423.2499

$ws.Range("H134").Value = 1460.1177
$ws.Range("I134").Value = 1412.7894
$ws.Range("J134").Value = 1520.0667
$ws.Range("K134").Value = 4238.3682
$ws.Range("L134").Value = 4560.2001
$ws.Range("M134").Value = -1703.3682
$ws.Range("N134").Value = -9630.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 350
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = -431
$ws.Range("N17").Value = -1838

$ws.Range("H34").Value = 4547257
$ws.Range("I34").Value = 347.125
$ws.Range("J34").Value = 7145491.5
$ws.Range("K34").Value = 1041.375
$ws.Range("L34").Value = 21436474.5
$ws.Range("M34").Value = -957.375
$ws.Range("N34").Value = -21436642.5

$ws.Range("H39").Value = 4200.25
$ws.Range("J39").Value = 4257.4287
$ws.Range("L39").Value = 12772.2861
$ws.Range("N39").Value = -13360.2861

$ws.Range("H55").Value = 3100
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3100
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 9300
$ws.Range("N55").Value = -9654
$ws.Range("M55").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8383.571
$ws.Range("I132").Value = 9963.866
$ws.Range("J132").Value = 4432.8335
$ws.Range("K132").Value = 29891.598
$ws.Range("L132").Value = 13298.5005
$ws.Range("M132").Value = -27361.598
$ws.Range("N132").Value = -18358.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1091.9231
$ws.Range("I93").Value = 1017.36365
$ws.Range("J93").Value = 1502
$ws.Range("K93").Value = 1017.36365
$ws.Range("L93").Value = 1502
$ws.Range("M93").Value = 230.63635
$ws.Range("N93").Value = -3998

$ws.Range("H122").Value = 70834584
$ws.Range("I122").Value = 94445110
$ws.Range("K122").Value = 283335330
$ws.Range("M122").Value = -283332880

$ws.Range("H132").Value = 2496.3076
$ws.Range("I132").Value = 2000.375
$ws.Range("J132").Value = 3289.8
$ws.Range("K132").Value = 6001.125
$ws.Range("L132").Value = 9869.400000000001
$ws.Range("M132").Value = -3471.125
$ws.Range("N132").Value = -14929.4

$ws.Range("H136").Value = 1746.1765
$ws.Range("I136").Value = 1607.8334
$ws.Range("J136").Value = 2078.2
$ws.Range("K136").Value = 4823.5002
$ws.Range("L136").Value = 6234.599999999999
$ws.Range("M136").Value = -2273.5002
$ws.Range("N136").Value = -11334.6

$ws.Range("H139").Value = 37447.855
$ws.Range("J139").Value = 37447.855
$ws.Range("L139").Value = 37447.855
$ws.Range("N139").Value = -47727.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10419097
$ws.Range("I122").Value = 13891390
$ws.Range("J122").Value = 2220
$ws.Range("K122").Value = 41674170
$ws.Range("L122").Value = 6660
$ws.Range("M122").Value = -41671720
$ws.Range("N122").Value = -11560

$ws.Range("H126").Value = 52911104
$ws.Range("I126").Value = 85470840
$ws.Range("K126").Value = 256412520
$ws.Range("M126").Value = -256410050

$ws.Range("H136").Value = 1106.9
$ws.Range("I136").Value = 853.5833
$ws.Range("J136").Value = 1486.875
$ws.Range("K136").Value = 2560.7499
$ws.Range("L136").Value = 4460.625
$ws.Range("M136").Value = -10.7498999999998
$ws.Range("N136").Value = -9560.625

$ws.Range("H141").Value = 54253.57
$ws.Range("J141").Value = 54253.57
$ws.Range("L141").Value = 54253.57
$ws.Range("N141").Value = -64613.57
